$wb = $excel.ActiveWorkbook

# Update the status text on every sheet where it appears ("Ready for handoff" -> "In Translation")
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Shrink the status columns to match the new (shorter) text width
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C:C").ColumnWidth = 12.5
